$d = $word.ActiveDocument

# The template began with 17 blank "filler" paragraphs (each just an
# empty w:p / w:pPr with szCs=22) followed by a final paragraph that
# holds the real tab stop used to lay out the form. The filler
# paragraphs are no longer needed, so remove all paragraphs up to (but
# not including) the last paragraph before the section break.

$paraCount = $d.Paragraphs.Count

if ($paraCount -gt 1) {
    $lastBodyParaIndex = $paraCount
    $startRange = $d.Paragraphs(1).Range.Start
    $endRange = $d.Paragraphs($lastBodyParaIndex - 1).Range.End
    $fillerRange = $d.Range($startRange, $endRange)
    $fillerRange.Delete()
}

# The remaining paragraph still contains a single tab character (a
# w:r with a w:tab) used with its tab stop; that run is no longer
# needed now that the preceding filler paragraphs are gone, so strip
# the tab character while keeping the paragraph (and its tab-stop
# formatting) intact.

$last = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $last.Range
$textLen = $lastRange.End - $lastRange.Start

if ($textLen -gt 1) {
    # Exclude the trailing paragraph mark from the deletion range.
    $tabOnly = $d.Range($lastRange.Start, $lastRange.End - 1)
    $tabOnly.Delete()
}
